$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $ok = $r.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found: [$findText]"
    }
    return $ok
}

# 1. "each." -> "each:" (split into 2 runs in the diff; content change period->colon)
Replace-Text "Which of the following statements are true/false, type T or F next to each." "Which of the following statements are true/false, type T or F next to each:" | Out-Null

# 2. Exercise 1 heading
Replace-Text "Exercise 1. Consequence of Openness" "Exercise 1. Benefits of openness" | Out-Null

# 3. Being open has other consequences...
Replace-Text "Being open has other consequences beyond giving the free access to information. " "Being open has other benefits beyond giving free access to information. " | Out-Null

# 4. For example "Open educational resources":
Replace-Text "For example “Open educational resources”:" "For example, “Open educational resources”:" | Out-Null

# 5. improves teachers/instructors
Replace-Text "- improves teachers/instructors skills by sharing ideas" "- improves teachers’/instructors’ skills by sharing ideas" | Out-Null

# 6. Discuss in your group
Replace-Text "Discuss in your group what are the additional benefits or addressed problems for the selected Open initiative:" "Discuss in your group what the additional benefits or addressed problems are for the selected open practices:" | Out-Null

# 7. Read through them
Replace-Text " Read through them, select 3 most important/attractive for you and mark them with +1, select two least important for you and mark them with 0" " Read through them, select the 3 most important/attractive for you and mark them with +1, select two least important for you and mark them with 0" | Out-Null

# 8. get extra value
Replace-Text "get extra value from your work (e.g. collaborators, reuse by modelers, ML specialists)" "get extra value from your work (e.g. collaborators, reuse by modellers, ML specialists)" | Out-Null

# 9. avoid embarassment
Replace-Text "avoid embarassment/disaster when you can’t reproduce your results" "avoid embarrassment/disaster when you cannot reproduce your results" | Out-Null

# 10. Can you think of other benefits
Replace-Text "Can you think of other benefits? How personal benefits of Open Science compare to the benefits for the (scientific) society." "Can you think of other benefits? How do personal benefits of Open Science compare to the benefits for the (scientific) society." | Out-Null

# 11. Open Science relies strongly on the Internet
Replace-Text "Open Science relies strongly on the Internet" "Open Science relies strongly on the internet" | Out-Null

# 12. You cannot Open Source patented software
Replace-Text "You cannot Open Source patented software" "You cannot Open-Source patented software" | Out-Null

# 13. You cannot charge for Open Source software
Replace-Text "You cannot charge for Open Source software" "You cannot charge for Open-Source software" | Out-Null

# 14. which one of the statements best characterize
Replace-Text "which one of the statements best characterize your experience (type +1" "which one of the statements best characterise your experience (type +1" | Out-Null
